$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New "type d'agent" labels (row 4 / row 5) ----------------------------
$ws.Range("B4").Value = "Référentiel dit des personnes morales"
$ws.Range("B5").Value = "Référentiel dit des personnes physiques"

# --- New "nombre d'entités" counts (row 4 / row 5) ------------------------
$ws.Range("E4").Value = 487
$ws.Range("E5").Value = 2134

# --- New descriptions (row 4 / row 5) -------------------------------------
$ws.Range("G4").Value = "Référentiel produit automatiquement à partir du RI_013 du SIA. Pas de différence notable avec le contenu du référentiel SIA. 40 entités alignées (owl:sameAs) avec autant d'entités du référentiel des producteurs."
$ws.Range("G5").Value = "Référentiel produit automatiquement à partir du RI_012 du SIA. Pas de différence notable avec le contenu du référentiel SIA. 102 entités alignées (owl:sameAs) avec autant d'entités du référentiel des producteurs."

# --- Row 6 (agents/producteurs) : fix typo + add "EAC-CPF" ---------------
$ws.Range("G6").Value = "Référentiel produit à partir du référentiel des producteurs du SIA, en utilisant RiC-O Converter, qui regroupe et déduplique les relations EAC-CPF. Un fichier par agent, et plusieurs fichiers pour les relations. Fourni avec une liste au format tsv (en utf-8; séparateur : tabulation) des producteurs"

# --- Row heights for the wrapped descriptions -----------------------------
$ws.Range("A4").EntireRow.RowHeight = 45
$ws.Range("A5").EntireRow.RowHeight = 45
$ws.Range("A6").EntireRow.RowHeight = 75

# --- Selection / view ------------------------------------------------------
$ws.Range("B8").Select()
